$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header labels (error message text) to lowercase
$ws.Range("B1").Value = "userid"
$ws.Range("C1").Value = "password"

# Update the selected cell to match the new active selection
$ws.Range("L8").Select()
